# Restore/apply change: cell C10 on sheet "Rules" should hold the numeric value 1
# (was previously 18).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

$ws.Range("C10").Value = 1
